$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.922.70"
$ws.Range("E2").Value = "  +0.10%  "

$ws.Range("D3").Value = "1.814.52"

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.17%  "

$ws.Range("E6").Value = "  +0.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4651"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3660"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.45%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07352"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.39%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8695"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.61%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.21"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.18%  "

$ws.Range("D12").Value = "1.818.83"
$ws.Range("E12").Value = "  -2.78%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.366"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.11%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07090"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.80%  "

$ws.Range("E15").Value = "  -0.10%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.67%  "

$ws.Range("E17").Value = "  +0.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008702"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.29%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.12%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.69%  "

$ws.Range("D21").Value = "26.931.59"
$ws.Range("E21").Value = "  +0.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.291"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.37%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.53%  "

$ws.Range("D24").Value = "2.024.03"
$ws.Range("E24").Value = "  -2.58%  "

$ws.Range("E25").Value = "  -0.68%  "

$ws.Range("E26").Value = "  -0.39%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.120"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.18%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.249"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.74%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.37%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08907"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.32%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7548"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.57%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.160"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.482"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.47%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.907"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.18%  "

$ws.Range("E36").Value = "  +0.21%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.087"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.79%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05284"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.49%  "

$ws.Range("E39").Value = "  -0.75%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.970"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.57%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.244"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5306"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.301"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.87%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1651"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.91%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.415"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.40%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4866"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.68%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.70%  "

$ws.Range("E48").Value = "  +0.20%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.658"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.47%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.25%  "

$ws.Range("E51").Value = "  -0.06%  "
